# Add a new last column (AA) for the date 12-10-2020, continuing the daily
# COVID-19 time series table. Column AA mirrors the formatting used by the
# rest of the header row (bold, thin border, centered alignment) and column
# Z in particular (the previous last date column), and rows 2-36 get the
# new per-state case counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID19_TIMESERIESDATA")

# --- Header cell AA1: the new date label, formatted like the rest of row 1 ---
$hdr = $ws.Range("AA1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.NumberFormat = "@"            # keep as text so it is not parsed as a date
$hdr.Value = "12-10-2020"

# --- Data rows 2-36: new case counts for 12-10-2020, one per State/UT ---
$values = @{
    2  = 186
    3  = 46295
    4  = 2891
    5  = 28385
    6  = 11044
    7  = 1184
    8  = 27348
    9  = 102
    10 = 21701
    11 = 4656
    12 = 15695
    13 = 10573
    14 = 2687
    15 = 10466
    16 = 8167
    17 = 120289
    18 = 96401
    19 = 980
    20 = 15177
    21 = 221637
    22 = 2731
    23 = 2478
    24 = 174
    25 = 1259
    26 = 23602
    27 = 4695
    28 = 9275
    29 = 21412
    30 = 384
    31 = 44095
    32 = 24514
    33 = 3742
    34 = 7373
    35 = 40019
    36 = 30236
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 27).Value = $values[$row]
}
